$wb = $excel.ActiveWorkbook

# --- Rename the sample user everywhere first ------------------------------
# "Melissa Zatta" becomes "Julie Carthane" on both sheets that reference her
# (and the cell loses its wrap-text style). Doing both occurrences before
# introducing any brand-new string keeps the renamed entry in its original
# shared-string slot instead of appending a fresh one.
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Julie Carthane"
$wsUsers.Range("A2").Style = "Normal"

$wsGiftLog = $wb.Worksheets.Item("GiftLog")
$wsGiftLog.Range("B2").Value = "Julie Carthane"
$wsGiftLog.Range("B2").Style = "Normal"
$wsGiftLog.Range("F20").Select() | Out-Null

# --- Users sheet header ----------------------------------------------------
# Header "StdUser" is replaced by "CF FinancialUser" (a brand-new string,
# appended after the rename above).
$wsUsers.Range("A1").Value = "CF FinancialUser"

# --- SearchCriteria sheet --------------------------------------------------
# This sheet was previously the active tab / selection; the new selection
# keeps its G6 cell (unchanged) but the tab focus moves away (handled below
# by activating the Users sheet last).
$wsSearchCriteria = $wb.Worksheets.Item("SearchCriteria")
$wsSearchCriteria.Range("G6").Select() | Out-Null

# --- Final focus -----------------------------------------------------------
# Users becomes the active / selected tab, with the cursor on A5.
$wsUsers.Activate() | Out-Null
$wsUsers.Range("A5").Select() | Out-Null
